$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.170.95'
$ws.Range('E2').Value = '  +2.39%  '

$ws.Range('D3').Value = '2.050.01'
$ws.Range('E3').Value = '  +1.67%  '

$ws.Range('E4').Value = '  +0.15%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '249.96'
$ws.Range('E5').Value = '  +0.21%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.644'
$ws.Range('E6').Value = '  -0.04%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '63.93'
$ws.Range('E7').Value = '  +1.01%  '

$ws.Range('E8').Value = '  +0.11%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.390'
$ws.Range('E9').Value = '  +5.29%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '58.75'
$ws.Range('E10').Value = '  -1.01%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0795'
$ws.Range('E11').Value = '  +6.13%  '

$ws.Range('E12').Value = '  -0.92%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.906'
$ws.Range('E13').Value = '  -4.48%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '22.81'
$ws.Range('E14').Value = '  +17.13%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.51'
$ws.Range('E15').Value = '  -3.87%  '

$ws.Range('D16').Value = '2.359.19'
$ws.Range('E16').Value = '  +2.19%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.55'
$ws.Range('E17').Value = '  +2.07%  '

$ws.Range('D18').Value = '2.057.32'
$ws.Range('E18').Value = '  +2.13%  '

$ws.Range('D19').Value = '37.132.08'
$ws.Range('E19').Value = '  +2.72%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.79'
$ws.Range('E20').Value = '  +0.85%  '

$ws.Range('D21').Value = '0.0₃0887'
$ws.Range('E21').Value = '  +3.19%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.42'
$ws.Range('E22').Value = '  +1.93%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.27'
$ws.Range('E23').Value = '  +1.03%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  -0.14%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.57'
$ws.Range('E25').Value = '  -4.33%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.36'
$ws.Range('E26').Value = '  +3.02%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.82'
$ws.Range('E27').Value = '  +1.51%  '

$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '159.88'
$ws.Range('E28').Value = '  -3.77%  '

$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.30'
$ws.Range('E29').Value = '  +3.13%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.131'
$ws.Range('E30').Value = '  +22.30%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.122'
$ws.Range('E31').Value = '  +0.99%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.08'
$ws.Range('E32').Value = '  -1.22%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.18'
$ws.Range('E33').Value = '  -1.23%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0621'
$ws.Range('E34').Value = '  +2.15%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.56'
$ws.Range('E35').Value = '  +1.46%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.39'
$ws.Range('E36').Value = '  -4.01%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.32'
$ws.Range('E37').Value = '  +8.46%  '

$ws.Range('E38').Value = '  +0.19%  '

$ws.Range('E39').Value = '  +1.82%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.95'
$ws.Range('E40').Value = '  +26.32%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.101'
$ws.Range('E41').Value = '  -3.81%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.25'
$ws.Range('E42').Value = '  +2.54%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.01'
$ws.Range('E43').Value = '  +3.81%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.15'
$ws.Range('E44').Value = '  +2.27%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '17.10'
$ws.Range('E45').Value = '  -0.78%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0216'
$ws.Range('E46').Value = '  -0.21%  '

$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '93.80'
$ws.Range('E47').Value = '  -0.98%  '

$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.78'
$ws.Range('E48').Value = '  -0.09%  '

$ws.Range('D49').Value = '1.378.87'
$ws.Range('E49').Value = '  -0.21%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.91'
$ws.Range('E50').Value = '  -0.93%  '

$ws.Range('D51').Value = '2.243.12'
$ws.Range('E51').Value = '  +2.16%  '
